$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Activate()

# Insert a new blank column before column N (14) - this shifts the
# existing "Late" / "Outstanding" columns one place to the right,
# matching the Variable-Instalments layout used elsewhere in the
# RBI loan automation workbooks.
$ws.Columns("N:N").Insert()

# The newly inserted column gets a wider, non-"best fit" width.
$ws.Columns("N:N").ColumnWidth = 9.1

# Re-assign the "Outstanding" values that were shifted from P into Q so
# they keep their original, clean decimal representation instead of
# picking up floating-point noise from the column move.
$ws.Range("Q3").Value = 935.25
$ws.Range("Q4").Value = 925.38
$ws.Range("Q5").Value = 935.25
$ws.Range("Q6").Value = 907.3
$ws.Range("Q7").Value = 901.28
$ws.Range("Q8").Value = 890.86
$ws.Range("Q9").Value = 884.29
$ws.Range("Q10").Value = 875.8
$ws.Range("Q11").Value = 866.21
$ws.Range("Q12").Value = 858.81
$ws.Range("Q13").Value = 849.77
$ws.Range("Q14").Value = 841.86

# Put the selection where the author left it after editing.
$ws.Range("S8").Select() | Out-Null
